$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("References")

$ws1.Range("E5").Value2 = 'Based on official disease reports to the WOAH'
$ws1.Range("E6").Value2 = 'CSF is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005.'
$ws1.Range("E7").Value2 = 'As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:'
$ws1.Range("E14").Value2 = 'Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently.'
$ws1.Range("E17").Value2 = 'For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}.'
$ws1.Range("E21").Value2 = 'A summary of the disease in animal hosts is given in the {ref008:WOAH Technical disease card}.'
$ws1.Range("E34").Value2 = 'Humans are not affected by this virus. Swine are the only species known to be susceptible. ({ref008:WOAH}).'
$ws1.Range("E42").Value2 = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the virus characteristics. '
$ws1.Range("E54").Value2 = 'Refer to the {ref008:WOAH Technical disease card} for a key summary of the disease transmission and epidemiological parameters.'
$ws1.Range("E66").Value2 = 'WOAH-prescribed tests for international trade include: the neutralising peroxidase-linked assay (NPLA), fluorescent antibody virus neutralisation (FAVN), and enzyme-linked immunosorbent assay (ELISA). CSF cross-reactive antibodies may be present in sera from pigs infected with other pestiviruses (eg BVDV, BDV) so positive screening results should be confirmed with a CSFV-specific test ({ref010:WOAH, Terrestrial Manual})'
$ws1.Range("E78").Value2 = 'Emergency vaccination against CSF can be used within the European Union for the control of the disease in the event of an outbreak. Modified-live vaccines (MLV) based on attenuated CSF virus strains are commericially available.  Marker vaccines that allow differentiation of infected from vaccinated animals (DIVA) are under development and one such vaccine, a chimeric pestivirus based on the E2 subunit of CSFV, has been licensed by the European Medicines Agency. Effective inactivated vaccines are not available ({ref010:WOAH, Terrestrial Manual}).  '
$ws1.Range("E92").Value2 = 'Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data.'
$ws2.Range("C2").Value2 = 'WOAH-WAHIS (WOAH World Animal Health Information System)'
$ws2.Range("C6").Value2 = 'WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$ws2.Range("C9").Value2 = 'WOAH (World Organisation for Animal Health) Technical Disease Card: Classical swine fever. 2021.'
$ws2.Range("C10").Value2 = 'WOAH (World Organisation for Animal Health), 2021. Classical swine fever. Chapter 15.2. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France'
$ws2.Range("C11").Value2 = 'WOAH (World Organisation for Animal Health), 2019. Classical swine fever. Chapter 3.08.03. WOAH Terrestrial Manual 2019. WOAH, Paris, France'
